# Delete the duplicate row 2973 (date 45442 appeared twice); this shifts
# all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2973).Delete()

# Update the view to match the post-edit selection / scroll position.
$ws.Range("A2973:XFD2973").Select()
$ws.Application.ActiveWindow.ScrollRow = 2947

$excel.ActiveWindow.WindowState = -4137
